$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.278.42'
$ws.Range("E2").Value = '  +3.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.635.45'
$ws.Range("E3").Value = '  +2.46%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.13'
$ws.Range("E5").Value = '  +2.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.27'
$ws.Range("E6").Value = '  +4.82%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.118'
$ws.Range("E9").Value = '  +8.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.398'
$ws.Range("E10").Value = '  +4.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").Value = '  +2.49%  '

$ws.Range("E12").Value = '  +1.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.93'
$ws.Range("E13").Value = '  +5.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000185'
$ws.Range("E14").Value = '  +19.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.106.37'
$ws.Range("E15").Value = '  +2.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.128.02'
$ws.Range("E16").Value = '  +3.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.559.04'
$ws.Range("E17").Value = '  -1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").Value = '  +2.34%  '

$ws.Range("E19").Value = '  +2.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.81'
$ws.Range("E20").Value = '  +2.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.25'
$ws.Range("E21").Value = '  +5.64%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.02'
$ws.Range("E23").Value = '  +0.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.71'
$ws.Range("E24").Value = '  +1.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.49'
$ws.Range("E25").Value = '  +3.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.65'
$ws.Range("E26").Value = '  -0.54%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  +1.83%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.09'
$ws.Range("E28").Value = '  +1.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0951'
$ws.Range("E29").Value = '  +12.99%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  +4.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '518.44'
$ws.Range("E32").Value = '  -5.85%  '

$ws.Range("E33").Value = '  +3.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.68'
$ws.Range("E34").Value = '  +9.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.35'
$ws.Range("E35").Value = '  +4.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.427'
$ws.Range("E36").Value = '  +4.32%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.05'
$ws.Range("E37").Value = '  +7.13%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.20'
$ws.Range("E38").Value = '  -1.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.21'
$ws.Range("E39").Value = '  +3.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.14'
$ws.Range("E42").Value = '  +6.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.19'
$ws.Range("E43").Value = '  -0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.09'
$ws.Range("E44").Value = '  +4.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0608'
$ws.Range("E45").Value = '  +4.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.12'
$ws.Range("E46").Value = '  +3.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.22'
$ws.Range("E47").Value = '  +9.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.648'
$ws.Range("E48").Value = '  +3.55%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0255'
$ws.Range("E49").Value = '  +1.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0983'
$ws.Range("E50").Value = '  +2.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.47'
$ws.Range("E51").Value = '  +2.76%  '
